$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1927194860813705
$ws.Range("C2").Value = 0.5674518201284796
$ws.Range("J2").Value = 0.01713062098501071
$ws.Range("P2").Value = 0.145610278372591
$ws.Range("S2").Value = 0.07708779443254818
$ws.Range("B3").Value = 0.007462686567164179
$ws.Range("C3").Value = 0.01865671641791045
$ws.Range("J3").Value = 0.007462686567164179
$ws.Range("P3").Value = 0.7910447761194029
$ws.Range("S3").Value = 0.1753731343283582
$ws.Range("J4").Value = 0.03225806451612903
$ws.Range("P4").Value = 0.8548387096774194
$ws.Range("S4").Value = 0.1129032258064516
$ws.Range("B6").Value = 0.0469208211143695
$ws.Range("D6").Value = 0.01466275659824047
$ws.Range("F6").Value = 0.06744868035190615
$ws.Range("J6").Value = 0.2521994134897361
$ws.Range("O6").Value = 0.01173020527859238
$ws.Range("Q6").Value = 0.1847507331378299
$ws.Range("R6").Value = 0.09970674486803519
$ws.Range("S6").Value = 0.3225806451612903
$ws.Range("B7").Value = 0.1150159744408946
$ws.Range("D7").Value = 0.01916932907348243
$ws.Range("F7").Value = 0.05111821086261981
$ws.Range("J7").Value = 0.134185303514377
$ws.Range("O7").Value = 0.01916932907348243
$ws.Range("Q7").Value = 0.2044728434504792
$ws.Range("R7").Value = 0.0670926517571885
$ws.Range("S7").Value = 0.389776357827476
$ws.Range("B8").Value = 0.1049723756906077
$ws.Range("D8").Value = 0.008287292817679558
$ws.Range("E8").Value = 0.002762430939226519
$ws.Range("F8").Value = 0.07458563535911603
$ws.Range("J8").Value = 0.1422651933701657
$ws.Range("O8").Value = 0.01104972375690608
$ws.Range("Q8").Value = 0.1947513812154696
$ws.Range("R8").Value = 0.1022099447513812
$ws.Range("S8").Value = 0.3591160220994475
$ws.Range("B9").Value = 0.1241830065359477
$ws.Range("D9").Value = 0.009803921568627451
$ws.Range("F9").Value = 0.05228758169934641
$ws.Range("J9").Value = 0.1176470588235294
$ws.Range("O9").Value = 0.0196078431372549
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.07516339869281045
$ws.Range("S9").Value = 0.434640522875817
$ws.Range("B10").Value = 0.1033099297893681
$ws.Range("D10").Value = 0.02156469408224674
$ws.Range("E10").Value = 0.0005015045135406219
$ws.Range("F10").Value = 0.07472417251755266
$ws.Range("J10").Value = 0.1038114343029087
$ws.Range("O10").Value = 0.01203610832497492
$ws.Range("Q10").Value = 0.2166499498495486
$ws.Range("R10").Value = 0.09327983951855567
$ws.Range("S10").Value = 0.3741223671013039
$ws.Range("G11").Value = 0.1340425531914894
$ws.Range("J11").Value = 0.09148936170212765
$ws.Range("K11").Value = 0.1638297872340425
$ws.Range("L11").Value = 0.6042553191489362
$ws.Range("S11").Value = 0.006382978723404255
$ws.Range("G12").Value = 0.7226027397260274
$ws.Range("J12").Value = 0.2157534246575342
$ws.Range("K12").Value = 0.00684931506849315
$ws.Range("L12").Value = 0.03082191780821918
$ws.Range("S12").Value = 0.02397260273972603
$ws.Range("G13").Value = 0.6575342465753424
$ws.Range("J13").Value = 0.3150684931506849
$ws.Range("S13").Value = 0.0273972602739726
$ws.Range("F15").Value = 0.01159420289855072
$ws.Range("H15").Value = 0.1768115942028985
$ws.Range("I15").Value = 0.06666666666666667
$ws.Range("J15").Value = 0.3594202898550725
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.01739130434782609
$ws.Range("N15").Value = 0.002898550724637681
$ws.Range("O15").Value = 0.06086956521739131
$ws.Range("S15").Value = 0.2376811594202899
$ws.Range("F16").Value = 0.01547987616099071
$ws.Range("H16").Value = 0.1671826625386997
$ws.Range("I16").Value = 0.09287925696594428
$ws.Range("J16").Value = 0.3653250773993808
$ws.Range("K16").Value = 0.1269349845201238
$ws.Range("M16").Value = 0.02476780185758514
$ws.Range("O16").Value = 0.04953560371517028
$ws.Range("S16").Value = 0.1578947368421053
$ws.Range("F17").Value = 0.01211305518169583
$ws.Range("H17").Value = 0.1830417227456258
$ws.Range("I17").Value = 0.08882907133243607
$ws.Range("J17").Value = 0.4078061911170929
$ws.Range("K17").Value = 0.1036339165545088
$ws.Range("M17").Value = 0.02153432032301481
$ws.Range("O17").Value = 0.06056527590847914
$ws.Range("S17").Value = 0.1224764468371467
$ws.Range("F18").Value = 0.04129793510324484
$ws.Range("H18").Value = 0.2005899705014749
$ws.Range("I18").Value = 0.0943952802359882
$ws.Range("J18").Value = 0.3864306784660767
$ws.Range("K18").Value = 0.08259587020648967
$ws.Range("M18").Value = 0.01474926253687316
$ws.Range("O18").Value = 0.07669616519174041
$ws.Range("S18").Value = 0.1032448377581121
$ws.Range("F19").Value = 0.01070336391437309
$ws.Range("H19").Value = 0.2089704383282365
$ws.Range("I19").Value = 0.07798165137614679
$ws.Range("J19").Value = 0.3700305810397553
$ws.Range("K19").Value = 0.109072375127421
$ws.Range("M19").Value = 0.02089704383282365
$ws.Range("O19").Value = 0.07186544342507645
$ws.Range("S19").Value = 0.1304791029561672
